# Daily attendance processing - 2025-10-31 19:42:15
#
# For every row in the "Recorded By" column (G), when the value is a
# comma-separated list of recorders that ends with the literal token
# "System", swap the last two entries - except for the specific
# "backup@backdoor.com, System" pairing, which is left untouched.
#
# Examples:
#   "dnasr281@gmail.com, System"              -> "System, dnasr281@gmail.com"
#   "admin@admin.com, System"                 -> "System, admin@admin.com"
#   "backup@backdoor.com, system, System"     -> "backup@backdoor.com, System, system"
#   "backup@backdoor.com, System"             -> (unchanged)
#   "dnasr281@gmail.com"                      -> (unchanged)
#   "System"                                  -> (unchanged)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $val = $ws.Cells.Item($row, 7).Value2

    if ($val -ne "" -and $val -ne $null) {
        $parts = $val.Split(",")

        if ($parts.Length -ge 2) {
            $lastPart = $parts[$parts.Length - 1].Trim()

            if (($lastPart -eq "System") -and ($val -ne "backup@backdoor.com, System")) {
                $newParts = @()
                for ($i = 0; $i -lt $parts.Length - 2; $i++) {
                    $newParts += $parts[$i].Trim()
                }
                $newParts += $parts[$parts.Length - 1].Trim()
                $newParts += $parts[$parts.Length - 2].Trim()

                $newVal = $newParts -join ", "
                $ws.Cells.Item($row, 7).Value = $newVal
            }
        }
    }
}
